# [AA | 14/05/2018] : Test Case Changes
#
# Adds a new "Process_FeaturesNode" worksheet in front of the existing
# "Process_SortNode" tab. It is populated with the same employee records,
# re-sorted by employee Id (column B, ascending), plus a new "feature"
# column G that is the concatenation of the employee Name (A) and
# Department (D) - e.g. "AMAN" & "Automation Anywhere" -> "AMANAutomation
# Anywhere". The new sheet becomes the active/selected tab.

$wb = $excel.ActiveWorkbook

# Insert the new sheet before the current first tab (Process_SortNode).
$new = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$new.Name = "Process_FeaturesNode"

# Name, Id, Department, Tool, Experience, Manager
# (leading "," on each row keeps the row as a nested array instead of
# being flattened into the outer @() list)
$data = @(
    ,@('AMAN',     100, 'RESEARCH',   'Automation Anywhere', '4 years',  'GIRISH')
    ,@('GIRISH',   101, 'SALES',      'UiPath',               '5 years',  'MARTIN')
    ,@('AMOL',     102, 'SALES',      'Automation Anywhere', '6 years',  'AMAN')
    ,@('SUYOG',    103, 'ACCOUNTING', 'Blue Prism',           '7 years',  'SANGEETA')
    ,@('AAVESH',   104, 'RESEARCH',   'Development',          '8 years',  'DAVID')
    ,@('SANTOSH',  105, 'SALES',      'Selenium',             '9 years',  'NONAME')
    ,@('SANGEETA', 106, 'RESEARCH',   'UiPath',               '10 years', 'SAURABH')
    ,@('NISHA',    107, 'ACCOUNTING', 'RPA',                  '11 years', 'TEJAS')
    ,@('MARTIN',   108, 'SALES',      'Automation Anywhere', '12 years', 'NOOO')
    ,@('MILLER',   109, 'ACCOUNTING', 'Blue Prism',           '13 years', 'MARTIN')
    ,@('SCOTT',    110, 'RESEARCH',   'UiPath',               '14 years', 'LOL')
    ,@('SACHIN',   111, 'ACCOUNTING', 'RPA',                  '20 years', 'MANU')
    ,@('ayusH',    112, 'SALES',      'Selenium',             '9 years',  'NONAA')
)

$r = 1
foreach ($row in $data) {
    $name = $row[0]
    $dept = $row[3]

    $new.Cells.Item($r, 1).Value = $name
    $new.Cells.Item($r, 2).Value = $row[1]
    $new.Cells.Item($r, 3).Value = $row[2]
    $new.Cells.Item($r, 4).Value = $dept
    $new.Cells.Item($r, 5).Value = $row[4]
    $new.Cells.Item($r, 6).Value = $row[5]
    $new.Cells.Item($r, 7).Value = "$name$dept"

    $r++
}

$new.Activate()
